# 6.1.1.2b — add a 2022 data column (S) to the table and restyle column R
# (years 2007-2021 -> now 2007-2022) so it matches the formatting already
# used by column P, the way Excel normalizes a freshly-typed column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-style column R (header/value/value rows) to match column P ---
#     (row4: years header, row5/row6: data rows) — P's format "wins" and R
#     is brought in line with it before the new column is appended.
$ws.Range("P4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null

$ws.Range("P5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null

$ws.Range("P6").Copy() | Out-Null
$ws.Range("R6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 2. Add the new column S, cloning the (now updated) column R formats ---
$ws.Range("R3").Copy() | Out-Null
$ws.Range("S3").PasteSpecial(-4122) | Out-Null

$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null

$ws.Range("R5").Copy() | Out-Null
$ws.Range("S5").PasteSpecial(-4122) | Out-Null

$ws.Range("R6").Copy() | Out-Null
$ws.Range("S6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 3. Fill in the new column's values ---
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 1.8
$ws.Range("S6").Value = 8.4

# --- 4. Match the saved selection state ---
$ws.Range("S3").Select() | Out-Null
